$wb = $excel.ActiveWorkbook

# BranchData: branch 3-4 (row 5) was removed from the table, so every row
# below it shifts up by one.
$ws1 = $wb.Worksheets.Item("BranchData")
$ws1.Rows("5:5").Delete()

# BusData tab loses focus/selection highlight (moves to K14) but stays
# inactive - set its selection first, before activating BranchData, so the
# final "active sheet" ends up being BranchData.
$ws2 = $wb.Worksheets.Item("BusData")
$ws2.Activate()
$ws2.Range("K14").Select()

# BranchData becomes the active/selected tab again, with the cursor on H19.
$ws1.Activate()
$ws1.Range("H19").Select()
